$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the species-record data (columns A, B, E, F, G, H, Q, R)
# between row 2 and row 3, and between row 4 and row 5, while leaving all
# other columns (C, D, I, J, K, N, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AF, AG, AT, AW, AX, AY) untouched.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")
$rowPairs = @(
    @(2, 3),
    @(4, 5)
)

foreach ($pair in $rowPairs) {
    $rowX = $pair[0]
    $rowY = $pair[1]

    foreach ($col in $cols) {
        $cellX = $ws.Range("$col$rowX")
        $cellY = $ws.Range("$col$rowY")

        $valX = $cellX.Value2
        $valY = $cellY.Value2

        $cellX.Value2 = $valY
        $cellY.Value2 = $valX
    }
}
